$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = "2022-09-02 21:00:48"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "6670192"
$ws.Range("B3").Value = "Philips Wasserkocher HD9318/01"
$ws.Range("C3").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/philips-wasserkocher-hd931801/p/6670192"
$ws.Range("G3").Value = "Philips"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "24.95"
$ws.Range("M3").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N3").Value = "Philips Wasserkocher HD9318/01 50% Aktion 24.95 Schweizer Franken statt 49.90 Schweizer Franken"
$ws.Range("O3").Value = "2022-09-02 21:00:48"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "6579165"
$ws.Range("B4").Value = "satrap Trimm Set XA 5-in-1 Haar/Trimm-Set"
$ws.Range("C4").Value = "/de/kosmetik-gesundheit/herrenpflege-rasur/elektrische-herrenrasierer/satrap-trimm-set-xa-5-in-1-haartrimm-set/p/6579165"
$ws.Range("E4").ClearContents() | Out-Null
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = "satrap"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "44.95"
$ws.Range("M4").Value = "['kosmetik-gesundheit', 'herrenpflege-rasur', 'elektrische-herrenrasierer']"
$ws.Range("N4").Value = "satrap Trimm Set XA 5-in-1 Haar/Trimm-Set 25% Aktion 44.95 Schweizer Franken statt 59.95 Schweizer Franken"
$ws.Range("O4").Value = "2022-09-02 21:00:48"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "6007534"
$ws.Range("B5").Value = "Trend USB-Stick 8 GB"
$ws.Range("C5").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-8-gb/p/6007534"
$ws.Range("D5").ClearContents() | Out-Null
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = "Trend"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "9.95"
$ws.Range("I5").ClearContents() | Out-Null
$ws.Range("J5").ClearContents() | Out-Null
$ws.Range("K5").ClearContents() | Out-Null
$ws.Range("L5").ClearContents() | Out-Null
$ws.Range("M5").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Range("N5").Value = "Trend USB-Stick 8 GB 9.95 Schweizer Franken"
$ws.Range("O5").Value = "2022-09-02 21:00:48"

# Row 6
$ws.Range("O6").Value = "2022-09-02 21:00:48"

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "3494230"
$ws.Range("B7").Value = "Varta Electronics V13GS / V357 1er Bli"
$ws.Range("C7").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v13gs--v357-1er-bli/p/3494230"
$ws.Range("D7").Value = "1ST"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "4.95"
$ws.Range("I7").Value = "4.95/1ST"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "4.95"
$ws.Range("M7").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N7").Value = "Varta Electronics V13GS / V357 1er Bli 4.95 Schweizer Franken"
$ws.Range("O7").Value = "2022-09-02 21:00:48"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "4589934"
$ws.Range("B8").Value = "Varta Longlife AA 4er Bli"
$ws.Range("C8").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-aa-4er-bli/p/4589934"
$ws.Range("D8").Value = "4ST"
$ws.Range("E8").ClearContents() | Out-Null
$ws.Range("F8").Value = 0
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "7.95"
$ws.Range("I8").Value = "1.99/1ST"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "1.99"
$ws.Range("M8").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N8").Value = "Varta Longlife AA 4er Bli 7.95 Schweizer Franken"
$ws.Range("O8").Value = "2022-09-02 21:00:48"

# Row 9
$ws.Range("O9").Value = "2022-09-02 21:00:48"

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "3591269"
$ws.Range("B10").Value = "Varta Longlife Max Power C 2er Bli"
$ws.Range("C10").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-max-power-c-2er-bli/p/3591269"
$ws.Range("D10").Value = "2ST"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = "Varta"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "8.95"
$ws.Range("I10").Value = "4.48/1ST"
$ws.Range("J10").Value = "Preis pro 1 Stück"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "4.48"
$ws.Range("L10").Value = "1ST"
$ws.Range("M10").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N10").Value = "Varta Longlife Max Power C 2er Bli 8.95 Schweizer Franken"
$ws.Range("O10").Value = "2022-09-02 21:00:48"

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "7016089"
$ws.Range("B11").Value = "Tefal Bügeleisen FV4961S0"
$ws.Range("C11").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/tefal-buegeleisen-fv4961s0/p/7016089"
$ws.Range("G11").Value = "Tefal"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "49.95"
$ws.Range("M11").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Range("N11").Value = "Tefal Bügeleisen FV4961S0 50% Aktion 49.95 Schweizer Franken statt 99.90 Schweizer Franken"
$ws.Range("O11").Value = "2022-09-02 21:00:48"

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "4358322"
$ws.Range("B12").Value = "Rayovac Hörgerätebatterien 13 6 Stück"
$ws.Range("C12").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/rayovac-hoergeraetebatterien-13-6-stueck/p/4358322"
$ws.Range("D12").Value = "6ST"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = "Rayovac"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "9.95"
$ws.Range("I12").Value = "1.66/1ST"
$ws.Range("J12").Value = "Preis pro 1 Stück"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = "1.66"
$ws.Range("L12").Value = "1ST"
$ws.Range("M12").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N12").Value = "Rayovac Hörgerätebatterien 13 6 Stück 9.95 Schweizer Franken"
$ws.Range("O12").Value = "2022-09-02 21:00:48"

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "6753975"
$ws.Range("B13").Value = "Severin Standgrill mit Grillplatte PG 8563"
$ws.Range("C13").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-standgrill-mit-grillplatte-pg-8563/p/6753975"
$ws.Range("D13").ClearContents() | Out-Null
$ws.Range("E13").ClearContents() | Out-Null
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = "Severin"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "74.50"
$ws.Range("I13").ClearContents() | Out-Null
$ws.Range("J13").ClearContents() | Out-Null
$ws.Range("K13").ClearContents() | Out-Null
$ws.Range("L13").ClearContents() | Out-Null
$ws.Range("M13").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N13").Value = "Severin Standgrill mit Grillplatte PG 8563 50% Aktion 74.50 Schweizer Franken statt 149.00 Schweizer Franken"
$ws.Range("O13").Value = "2022-09-02 21:00:48"

# Row 14
$ws.Range("O14").Value = "2022-09-02 21:00:48"

# Row 15
$ws.Range("O15").Value = "2022-09-02 21:00:48"

# Row 16
$ws.Range("O16").Value = "2022-09-02 21:00:48"

# Row 17
$ws.Range("O17").Value = "2022-09-02 21:00:48"

# Row 18
$ws.Range("O18").Value = "2022-09-02 21:00:48"

# Row 19
$ws.Range("O19").Value = "2022-09-02 21:00:48"

# Row 20
$ws.Range("O20").Value = "2022-09-02 21:00:48"

# Row 21
$ws.Range("O21").Value = "2022-09-02 21:00:48"

# Row 22
$ws.Range("O22").Value = "2022-09-02 21:00:48"

# Row 23
$ws.Range("O23").Value = "2022-09-02 21:00:48"

# Row 24
$ws.Range("O24").Value = "2022-09-02 21:00:48"

# Row 25
$ws.Range("O25").Value = "2022-09-02 21:00:48"

# Row 26
$ws.Range("O26").Value = "2022-09-02 21:00:48"

# Row 27
$ws.Range("O27").Value = "2022-09-02 21:00:48"

# Row 28
$ws.Range("O28").Value = "2022-09-02 21:00:48"

# Row 29
$ws.Range("O29").Value = "2022-09-02 21:00:48"

# Row 30
$ws.Range("O30").Value = "2022-09-02 21:00:48"

# Row 31
$ws.Range("O31").Value = "2022-09-02 21:00:48"

# Row 32
$ws.Range("O32").Value = "2022-09-02 21:00:48"

# Row 33
$ws.Range("O33").Value = "2022-09-02 21:00:48"

# Row 34
$ws.Range("O34").Value = "2022-09-02 21:00:48"

# Row 35
$ws.Range("O35").Value = "2022-09-02 21:00:48"

# Row 36
$ws.Range("O36").Value = "2022-09-02 21:00:48"
